$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.111.51"
$ws.Range("E2").Value = "  -4.71%  "
$ws.Range("D3").Value = "1.656.63"
$ws.Range("E3").Value = "  -3.13%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "'218.16"
$ws.Range("E5").Value = "  -2.87%  "
$ws.Range("D6").Value = "'0.5163"
$ws.Range("E6").Value = "  -3.10%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "'0.2574"
$ws.Range("E8").Value = "  -3.56%  "
$ws.Range("D9").Value = "'0.06432"
$ws.Range("E9").Value = "  -2.95%  "
$ws.Range("E10").Value = "  -5.04%  "
$ws.Range("D11").Value = "'0.07758"
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").Value = "1.656.21"
$ws.Range("E12").Value = "  -3.16%  "
$ws.Range("D13").Value = "'4.298"
$ws.Range("E13").Value = "  -5.58%  "
$ws.Range("D14").Value = "1.883.30"
$ws.Range("E14").Value = "  -3.19%  "
$ws.Range("D15").Value = "'0.5537"
$ws.Range("E15").Value = "  -4.18%  "
$ws.Range("D16").Value = "0.0₅8043"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "'64.31"
$ws.Range("E17").Value = "  -5.25%  "
$ws.Range("D18").Value = "26.161.77"
$ws.Range("E18").Value = "  -4.51%  "
$ws.Range("D19").Value = "'1.007"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "'210.65"
$ws.Range("E20").Value = "  -3.63%  "
$ws.Range("D21").Value = "'4.393"
$ws.Range("E21").Value = "  -5.62%  "
$ws.Range("D22").Value = "'10.05"
$ws.Range("E22").Value = "  -3.83%  "
$ws.Range("D23").Value = "'5.916"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "'143.87"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").Value = "'1.760"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("D27").Value = "'0.1159"
$ws.Range("E27").Value = "  -4.26%  "
$ws.Range("D28").Value = "'6.974"
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("D29").Value = "'15.74"
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("D30").Value = "'0.05287"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("D31").Value = "'1.252"
$ws.Range("E31").Value = "  -3.14%  "
$ws.Range("D32").Value = "'3.358"
$ws.Range("D33").Value = "'3.237"
$ws.Range("E33").Value = "  -5.56%  "
$ws.Range("D34").Value = "'1.574"
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("D35").Value = "'2.758"
$ws.Range("E35").Value = "  -4.14%  "
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("D37").Value = "'0.9230"
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.163.61"
$ws.Range("E38").Value = "  +10.96%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'0.5681"
$ws.Range("E39").Value = "  -3.34%  "
$ws.Range("D40").Value = "'0.01590"
$ws.Range("E40").Value = "  -2.88%  "
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").Value = "'0.8377"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").Value = "'5.647"
$ws.Range("E43").Value = "  -3.54%  "
$ws.Range("D44").Value = "'99.88"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").Value = "1.794.18"
$ws.Range("E45").Value = "  -3.19%  "
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  -7.39%  "
$ws.Range("D47").Value = "'0.4511"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").Value = "'55.95"
$ws.Range("E48").Value = "  -3.53%  "
$ws.Range("D49").Value = "'1.007"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "'7.876"
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("D51").Value = "'0.05068"
$ws.Range("E51").Value = "  -3.01%  "
